$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("posts")

# Update the COLOR column (F) values for three posts with new hex colors
# (order matters so new shared-string entries are appended in the same
# sequence as the target workbook)
$ws.Range("F6").Value = "#C7CCEC"
$ws.Range("F2").Value = "#C2D7CE"
$ws.Range("F12").Value = "#BBE6DD"

# Match the resulting selection state on the "posts" sheet
$ws.Range("E13").Select()
